$d = $word.ActiveDocument

# Word's wmlNamespace, needed on every fragment we hand to InsertXML since
# each call is parsed as a standalone XML fragment.
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. "Kryssreferanser er som regel ..." paragraph -----------------------
# Re-split the two original runs at "arkivdeler" and "tilstede" and wrap
# those words with proofErr (spelling / grammar) markers, as a spell/grammar
# pass over the paragraph would produce. The visible text is unchanged.
$targetText = "Kryssreferanser er som regel"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*$targetText*") {
        $xml = '<w:p ' + $wNs + '>' + `
            '<w:r><w:t xml:space="preserve">Kryssreferanser er som regel fra en mappe til en annen, men kan skje til/fra </w:t></w:r>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:t>arkivdeler</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:t xml:space="preserve"> eller registreringer også</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve">. Hvis alle kryssreferansene er </w:t></w:r>' + `
            '<w:proofErr w:type="gramStart"/>' + `
            '<w:r><w:t>tilstede</w:t></w:r>' + `
            '<w:proofErr w:type="gramEnd"/>' + `
            '<w:r><w:t xml:space="preserve"> i uttrekket er alt her greit, ellers</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve"> bør det komme varsel og manuell sjekk av hva det refereres til. </w:t></w:r>' + `
            '</w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

# --- 2. "BaseX A_Kr2_Eksisterer_kryssreferansene_i_uttrekket.xq" paragraph -
# Same treatment: split at "BaseX" and "uttrekket.xq" and wrap with proofErr
# markers. Visible text is unchanged.
$targetText2 = "A_Kr2_Eksisterer_kryssreferansene"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*$targetText2*") {
        $xml = '<w:p ' + $wNs + '>' + `
            '<w:proofErr w:type="spellStart"/>' + `
            '<w:r><w:t>BaseX</w:t></w:r>' + `
            '<w:proofErr w:type="spellEnd"/>' + `
            '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
            '<w:r><w:t>A_Kr2_Eksisterer_kryssreferansene_i_</w:t></w:r>' + `
            '<w:proofErr w:type="gramStart"/>' + `
            '<w:r><w:t>uttrekket.xq</w:t></w:r>' + `
            '<w:proofErr w:type="gramEnd"/>' + `
            '<w:r><w:t xml:space="preserve"> for sjekking av kryssreferanser. Hvis det er for mange som mangler i uttrekket (f.eks. over 25) bør disse skrives ut som vedlegg i stedet. </w:t></w:r>' + `
            '</w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

# --- 3. Remove the stray "AND/OR" run left over in the last list paragraph -
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*AND/OR*") {
        $xml = '<w:p ' + $wNs + '><w:pPr><w:ind w:left="278" w:hanging="278"/></w:pPr></w:p>'
        $p.Range.InsertXML($xml)
        break
    }
}

Write-Output "edit applied"
